# Update "想去人数" (people-interested count) values in column F of the
# two sheets that carry the exhibition data: "展览" (sheet 1) and
# "全部类型" (sheet 4). Sheets "演出" and "本地生活" only contain a header
# row, so nothing to change there.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value, applied identically on both sheets.
$updates = @{
    4  = 1030
    5  = 27
    7  = 2671
    9  = 1694
    11 = 71
    12 = 567
    14 = 12
    15 = 73
    16 = 75
    17 = 80
}

foreach ($sheetIndex in 1, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
